# Auto-generated edit script applying numeric corrections to Lamia_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 589.6
$ws.Cells.Item(2, 9).Value = 112
$ws.Cells.Item(2, 11).Value = 112
$ws.Cells.Item(2, 13).Value = 1
$ws.Cells.Item(5, 8).Value = 34.625
$ws.Cells.Item(5, 9).Value = 34.625
$ws.Cells.Item(5, 11).Value = 34.625
$ws.Cells.Item(5, 13).Value = 80.375
$ws.Cells.Item(33, 8).Value = 449.46667
$ws.Cells.Item(33, 9).Value = 202.53847
$ws.Cells.Item(33, 11).Value = 202.53847
$ws.Cells.Item(33, 13).Value = 26.46153000000001
$ws.Cells.Item(86, 8).Value = 5873.231
$ws.Cells.Item(86, 10).Value = 5094
$ws.Cells.Item(86, 12).Value = 5094
$ws.Cells.Item(86, 14).Value = -7340
$ws.Cells.Item(89, 8).Value = 5873.231
$ws.Cells.Item(89, 10).Value = 5094
$ws.Cells.Item(89, 12).Value = 25470
$ws.Cells.Item(89, 14).Value = -36702
$ws.Cells.Item(92, 8).Value = 2858.5715
$ws.Cells.Item(92, 9).Value = 346
$ws.Cells.Item(92, 10).Value = 9140
$ws.Cells.Item(92, 11).Value = 346
$ws.Cells.Item(92, 12).Value = 9140
$ws.Cells.Item(92, 13).Value = 902
$ws.Cells.Item(92, 14).Value = -11636
$ws.Cells.Item(132, 8).Value = 1800.6897
$ws.Cells.Item(132, 9).Value = 1444.1538
$ws.Cells.Item(132, 11).Value = 4332.4614
$ws.Cells.Item(132, 13).Value = -1802.4614
$ws.Cells.Item(137, 8).Value = 3517.6086
$ws.Cells.Item(137, 9).Value = 4658
$ws.Cells.Item(137, 11).Value = 13974
$ws.Cells.Item(137, 13).Value = -11424
$ws.Cells.Item(141, 8).Value = 4668.0605
$ws.Cells.Item(141, 9).Value = 2133.353
$ws.Cells.Item(141, 10).Value = 7361.1875
$ws.Cells.Item(141, 11).Value = 6400.059
$ws.Cells.Item(141, 12).Value = 22083.5625
$ws.Cells.Item(141, 13).Value = -1220.059
$ws.Cells.Item(141, 14).Value = -32443.5625
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(13, 8).Value = 7502
$ws.Cells.Item(13, 10).Value = 10004
$ws.Cells.Item(13, 12).Value = 10004
$ws.Cells.Item(13, 14).Value = -10292
$ws.Cells.Item(25, 8).Value = 2271.6667
$ws.Cells.Item(25, 10).Value = 3000
$ws.Cells.Item(25, 12).Value = 3000
$ws.Cells.Item(25, 14).Value = -3804
$ws.Cells.Item(31, 8).Value = 3258.5
$ws.Cells.Item(31, 9).Value = 3258.5
$ws.Cells.Item(31, 11).Value = 3258.5
$ws.Cells.Item(31, 13).Value = -2964.5
$ws.Cells.Item(32, 8).Value = 2137.9666
$ws.Cells.Item(32, 9).Value = 1436.1177
$ws.Cells.Item(32, 10).Value = 6115.1113
$ws.Cells.Item(32, 11).Value = 1436.1177
$ws.Cells.Item(32, 12).Value = 6115.1113
$ws.Cells.Item(32, 13).Value = -1149.1177
$ws.Cells.Item(32, 14).Value = -6689.1113
$ws.Cells.Item(35, 8).Value = 8672.666999999999
$ws.Cells.Item(35, 9).Value = 2009
$ws.Cells.Item(35, 11).Value = 2009
$ws.Cells.Item(35, 13).Value = -1603
$ws.Cells.Item(61, 8).Value = 3679.12
$ws.Cells.Item(61, 9).Value = 3053.4783
$ws.Cells.Item(61, 11).Value = 3053.4783
$ws.Cells.Item(61, 13).Value = -2841.4783
$ws.Cells.Item(132, 8).Value = 2428.6667
$ws.Cells.Item(132, 9).Value = 1400.7142
$ws.Cells.Item(132, 10).Value = 6026.5
$ws.Cells.Item(132, 11).Value = 4202.142599999999
$ws.Cells.Item(132, 12).Value = 18079.5
$ws.Cells.Item(132, 13).Value = -1672.142599999999
$ws.Cells.Item(132, 14).Value = -23139.5
$ws.Cells.Item(136, 8).Value = 3679.12
$ws.Cells.Item(136, 9).Value = 3053.4783
$ws.Cells.Item(136, 11).Value = 9160.4349
$ws.Cells.Item(136, 13).Value = -6610.4349
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1301.3684
$ws.Cells.Item(99, 9).Value = 982.6
$ws.Cells.Item(99, 10).Value = 1655.5555
$ws.Cells.Item(99, 11).Value = 982.6
$ws.Cells.Item(99, 12).Value = 1655.5555
$ws.Cells.Item(99, 13).Value = 515.4
$ws.Cells.Item(99, 14).Value = -4651.5555
$ws.Cells.Item(107, 8).Value = 619.44446
$ws.Cells.Item(107, 9).Value = 584.375
$ws.Cells.Item(107, 10).Value = 900
$ws.Cells.Item(107, 11).Value = 584.375
$ws.Cells.Item(107, 12).Value = 900
$ws.Cells.Item(107, 13).Value = 1335.625
$ws.Cells.Item(107, 14).Value = -4740
$ws.Cells.Item(134, 8).Value = 3539.2917
$ws.Cells.Item(134, 9).Value = 1895.85
$ws.Cells.Item(134, 10).Value = 11756.5
$ws.Cells.Item(134, 11).Value = 5687.549999999999
$ws.Cells.Item(134, 12).Value = 35269.5
$ws.Cells.Item(134, 13).Value = -3152.549999999999
$ws.Cells.Item(134, 14).Value = -40339.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 7266
$ws.Cells.Item(22, 9).Value = 7499
$ws.Cells.Item(22, 11).Value = 7499
$ws.Cells.Item(22, 13).Value = -7149
$ws.Cells.Item(50, 8).Value = 26419
$ws.Cells.Item(50, 10).Value = 31987.5
$ws.Cells.Item(50, 12).Value = 31987.5
$ws.Cells.Item(50, 14).Value = -33237.5
$ws.Cells.Item(58, 8).Value = 9562.333000000001
$ws.Cells.Item(58, 9).Value = 2371
$ws.Cells.Item(58, 10).Value = 14699
$ws.Cells.Item(58, 11).Value = 2371
$ws.Cells.Item(58, 12).Value = 14699
$ws.Cells.Item(58, 13).Value = -2168
$ws.Cells.Item(58, 14).Value = -15105
$ws.Cells.Item(70, 8).Value = 44544.5
$ws.Cells.Item(70, 10).Value = 44544.5
$ws.Cells.Item(70, 12).Value = 44544.5
$ws.Cells.Item(70, 14).Value = -45174.5
$ws.Cells.Item(73, 8).Value = 44544.5
$ws.Cells.Item(73, 10).Value = 44544.5
$ws.Cells.Item(73, 12).Value = 44544.5
$ws.Cells.Item(73, 14).Value = -46728.5
$ws.Cells.Item(134, 8).Value = 2905.04
$ws.Cells.Item(134, 9).Value = 1914.4286
$ws.Cells.Item(134, 11).Value = 5743.2858
$ws.Cells.Item(134, 13).Value = -3208.2858
$ws.Cells.Item(136, 8).Value = 9562.333000000001
$ws.Cells.Item(136, 9).Value = 2371
$ws.Cells.Item(136, 10).Value = 14699
$ws.Cells.Item(136, 11).Value = 7113
$ws.Cells.Item(136, 12).Value = 44097
$ws.Cells.Item(136, 13).Value = -4563
$ws.Cells.Item(136, 14).Value = -49197
$ws.Cells.Item(141, 8).Value = 234798.8
$ws.Cells.Item(141, 10).Value = 234798.8
$ws.Cells.Item(141, 12).Value = 234798.8
$ws.Cells.Item(141, 14).Value = -245158.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 10926.5
$ws.Cells.Item(5, 10).Value = 41670.332
$ws.Cells.Item(5, 12).Value = 125010.996
$ws.Cells.Item(5, 14).Value = -125234.996
$ws.Cells.Item(122, 8).Value = 1642.2307
$ws.Cells.Item(122, 10).Value = 1712.7084
$ws.Cells.Item(122, 12).Value = 15414.3756
$ws.Cells.Item(122, 14).Value = -20314.3756
$ws.Cells.Item(135, 8).Value = 10926.5
$ws.Cells.Item(135, 10).Value = 41670.332
$ws.Cells.Item(135, 12).Value = 375032.988
$ws.Cells.Item(135, 14).Value = -380102.988
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 15128.8
$ws.Cells.Item(70, 9).Value = 5184.2144
$ws.Cells.Item(70, 11).Value = 5184.2144
$ws.Cells.Item(70, 13).Value = -4914.2144
$ws.Cells.Item(73, 8).Value = 15128.8
$ws.Cells.Item(73, 9).Value = 5184.2144
$ws.Cells.Item(73, 11).Value = 5184.2144
$ws.Cells.Item(73, 13).Value = -4248.2144
$ws.Cells.Item(97, 8).Value = 1455.4231
$ws.Cells.Item(97, 9).Value = 1148.421
$ws.Cells.Item(97, 11).Value = 1148.421
$ws.Cells.Item(97, 13).Value = -652.421
$ws.Cells.Item(107, 8).Value = 959.5
$ws.Cells.Item(107, 9).Value = 385.66666
$ws.Cells.Item(107, 10).Value = 1648.1
$ws.Cells.Item(107, 11).Value = 385.66666
$ws.Cells.Item(107, 12).Value = 1648.1
$ws.Cells.Item(107, 13).Value = 1534.33334
$ws.Cells.Item(107, 14).Value = -5488.1
$ws.Cells.Item(126, 8).Value = 4425
$ws.Cells.Item(126, 9).Value = 1914.5
$ws.Cells.Item(126, 10).Value = 8190.75
$ws.Cells.Item(126, 11).Value = 5743.5
$ws.Cells.Item(126, 12).Value = 24572.25
$ws.Cells.Item(126, 13).Value = -3273.5
$ws.Cells.Item(126, 14).Value = -29512.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5208.32
$ws.Cells.Item(7, 9).Value = 1749.8
$ws.Cells.Item(7, 10).Value = 19042.4
$ws.Cells.Item(7, 11).Value = 1749.8
$ws.Cells.Item(7, 12).Value = 19042.4
$ws.Cells.Item(7, 13).Value = -1637.8
$ws.Cells.Item(7, 14).Value = -19266.4
$ws.Cells.Item(16, 8).Value = 3173.9167
$ws.Cells.Item(16, 9).Value = 2209.889
$ws.Cells.Item(16, 11).Value = 2209.889
$ws.Cells.Item(16, 13).Value = -2039.889
$ws.Cells.Item(46, 8).Value = 5310.1
$ws.Cells.Item(46, 9).Value = 3499.75
$ws.Cells.Item(46, 10).Value = 6517
$ws.Cells.Item(46, 11).Value = 3499.75
$ws.Cells.Item(46, 12).Value = 6517
$ws.Cells.Item(46, 13).Value = -3311.75
$ws.Cells.Item(46, 14).Value = -6893
$ws.Cells.Item(61, 8).Value = 7992.273
$ws.Cells.Item(61, 9).Value = 6294.625
$ws.Cells.Item(61, 11).Value = 6294.625
$ws.Cells.Item(61, 13).Value = -6092.625
$ws.Cells.Item(113, 8).Value = 7992.273
$ws.Cells.Item(113, 9).Value = 6294.625
$ws.Cells.Item(113, 11).Value = 6294.625
$ws.Cells.Item(113, 13).Value = -4124.625
$ws.Cells.Item(126, 8).Value = 5208.32
$ws.Cells.Item(126, 9).Value = 1749.8
$ws.Cells.Item(126, 10).Value = 19042.4
$ws.Cells.Item(126, 11).Value = 5249.4
$ws.Cells.Item(126, 12).Value = 57127.2
$ws.Cells.Item(126, 13).Value = -2779.4
$ws.Cells.Item(126, 14).Value = -62067.2
$ws.Cells.Item(132, 8).Value = 5117
$ws.Cells.Item(132, 9).Value = 1944.5
$ws.Cells.Item(132, 10).Value = 8742.714
$ws.Cells.Item(132, 11).Value = 5833.5
$ws.Cells.Item(132, 12).Value = 26228.142
$ws.Cells.Item(132, 13).Value = -3303.5
$ws.Cells.Item(132, 14).Value = -31288.142
$ws.Cells.Item(135, 8).Value = 60301.875
$ws.Cells.Item(135, 10).Value = 60301.875
$ws.Cells.Item(135, 12).Value = 60301.875
$ws.Cells.Item(135, 14).Value = -70441.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 11300
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 11300
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 11300
$ws.Cells.Item(14, 13).ClearContents()
$ws.Cells.Item(14, 14).Value = -11636
$ws.Cells.Item(92, 8).Value = 30000
$ws.Cells.Item(92, 10).Value = 30000
$ws.Cells.Item(92, 12).Value = 30000
$ws.Cells.Item(92, 14).Value = -34992
$ws.Cells.Item(96, 8).Value = 974.8333
$ws.Cells.Item(96, 9).Value = 569.8
$ws.Cells.Item(96, 11).Value = 569.8
$ws.Cells.Item(96, 13).Value = 803.2
$ws.Cells.Item(122, 8).Value = 11746.6
$ws.Cells.Item(122, 10).Value = 17000.75
$ws.Cells.Item(122, 12).Value = 51002.25
$ws.Cells.Item(122, 14).Value = -55902.25
$ws.Cells.Item(132, 8).Value = 4526.4814
$ws.Cells.Item(132, 9).Value = 4579.7827
$ws.Cells.Item(132, 10).Value = 4220
$ws.Cells.Item(132, 11).Value = 13739.3481
$ws.Cells.Item(132, 12).Value = 12660
$ws.Cells.Item(132, 13).Value = -11209.3481
$ws.Cells.Item(132, 14).Value = -17720

Write-Host "Applied all corrections"